$wb = $excel.ActiveWorkbook

# Sheet "screwed" - append two new rows (index 9 -> stav, index 10 -> stav)
$wsScrewed = $wb.Worksheets.Item("screwed")
$wsScrewed.Range("A10").Value = 9
$wsScrewed.Range("B10").Value = "stav"
$wsScrewed.Range("A11").Value = 10
$wsScrewed.Range("B11").Value = "stav"

# Sheet "shiftsPerWorker" - update shift counts per worker
$wsShifts = $wb.Worksheets.Item("shiftsPerWorker")
$wsShifts.Range("B2").Value = 6
$wsShifts.Range("B3").Value = 4
$wsShifts.Range("B4").Value = 4
$wsShifts.Range("B5").Value = 4
$wsShifts.Range("B6").Value = 5
$wsShifts.Range("B7").Value = 1
